# "feat: add 2022-Q4 data"
#
# Original workbook sheets: 总计 (sheetId 1), 2022-Q2 (sheetId 2).
# Target workbook sheets:   总计 (sheetId 1), 2022-Q4 (sheetId 2), 2022-Q2 (sheetId 3).
#
# The existing "2022-Q2" sheet is renamed to "2022-Q4" and repopulated with
# the new quarter's fund-holdings data; a fresh sheet named "2022-Q2" is
# appended right after it, holding the fund-holdings data that used to live
# in the original "2022-Q2" sheet. The "总计" summary sheet gets a new row
# for 2022-Q4 inserted right after the header (pushing the existing
# 2022-Q2 row down one row).

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item(1)
$q2 = $wb.Worksheets.Item(2)

# Stash the original "2022-Q2" header/index-column formatting in a scratch
# range so it can be re-applied to the brand-new "2022-Q2" sheet created
# further down (new sheets don't inherit the old formatting automatically).
$q2.Range("B1:H1").Copy()
$q2.Range("Z1:AF1").PasteSpecial(-4122)
$q2.Cells.Item(2, 1).Copy()
$q2.Range("Z2").PasteSpecial(-4122)

# Rename the existing "2022-Q2" sheet to "2022-Q4" and load the new data.
$q2.Name = "2022-Q4"
$q4 = $q2

# The "2022-Q4" sheet picks up the "总计" sheet's header/index-column style.
$summary.Range("B1:D1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$summary.Cells.Item(2, 1).Copy()
$q4.Range("A2:A5").PasteSpecial(-4122)

$q4.Cells.Item(1, 2).Value = "基金代码"
$q4.Cells.Item(1, 3).Value = "基金名称"
$q4.Cells.Item(1, 4).Value = "基金规模"
$q4.Cells.Item(1, 5).Value = "股票总仓位"
$q4.Cells.Item(1, 6).Value = "仓位占比"
$q4.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q4.Cells.Item(1, 8).Value = "仓位排名"

$q4.Range("B2:B5").NumberFormat = "@"
$q4.Range("D2:G5").NumberFormat = "@"

$q4.Cells.Item(2, 1).Value = 0
$q4.Cells.Item(2, 2).Value = "213003"
$q4.Cells.Item(2, 3).Value = "宝盈策略增长混合"
$q4.Cells.Item(2, 4).Value = "10.36"
$q4.Cells.Item(2, 5).Value = "91.66"
$q4.Cells.Item(2, 6).Value = "4.10"
$q4.Cells.Item(2, 7).Value = "0.4248"
$q4.Cells.Item(2, 8).Value = 9

$q4.Cells.Item(3, 1).Value = 1
$q4.Cells.Item(3, 2).Value = "213002"
$q4.Cells.Item(3, 3).Value = "宝盈泛沿海增长混合"
$q4.Cells.Item(3, 4).Value = "5.08"
$q4.Cells.Item(3, 5).Value = "92.58"
$q4.Cells.Item(3, 6).Value = "4.89"
$q4.Cells.Item(3, 7).Value = "0.2484"
$q4.Cells.Item(3, 8).Value = 8

$q4.Cells.Item(4, 1).Value = 2
$q4.Cells.Item(4, 2).Value = "000796"
$q4.Cells.Item(4, 3).Value = "宝盈睿丰创新灵活配置混合 - C"
$q4.Cells.Item(4, 4).Value = "0.61"
$q4.Cells.Item(4, 5).Value = "92.26"
$q4.Cells.Item(4, 6).Value = "4.88"
$q4.Cells.Item(4, 7).Value = "0.0298"
$q4.Cells.Item(4, 8).Value = 8

$q4.Cells.Item(5, 1).Value = 3
$q4.Cells.Item(5, 2).Value = "000794"
$q4.Cells.Item(5, 3).Value = "宝盈睿丰创新灵活配置混合 - A/B"
$q4.Cells.Item(5, 4).Value = "0.39"
$q4.Cells.Item(5, 5).Value = "92.26"
$q4.Cells.Item(5, 6).Value = "4.88"
$q4.Cells.Item(5, 7).Value = "0.0190"
$q4.Cells.Item(5, 8).Value = 8

# Add a fresh "2022-Q2" sheet right after "2022-Q4" with the fund-holdings
# data that used to be in the original "2022-Q2" sheet.
$newQ2 = $wb.Worksheets.Add($null, $q4)
$newQ2.Name = "2022-Q2"

$q4.Range("Z1:AF1").Copy()
$newQ2.Range("B1:H1").PasteSpecial(-4122)
$q4.Range("Z2").Copy()
$newQ2.Cells.Item(2, 1).PasteSpecial(-4122)

$q4.Range("Z1:AF2").Clear()

$newQ2.Cells.Item(1, 2).Value = "基金代码"
$newQ2.Cells.Item(1, 3).Value = "基金名称"
$newQ2.Cells.Item(1, 4).Value = "基金规模"
$newQ2.Cells.Item(1, 5).Value = "股票总仓位"
$newQ2.Cells.Item(1, 6).Value = "仓位占比"
$newQ2.Cells.Item(1, 7).Value = "持有市值(亿元)"
$newQ2.Cells.Item(1, 8).Value = "仓位排名"

$newQ2.Range("B2").NumberFormat = "@"
$newQ2.Range("D2:G2").NumberFormat = "@"

$newQ2.Cells.Item(2, 1).Value = 0
$newQ2.Cells.Item(2, 2).Value = "501076"
$newQ2.Cells.Item(2, 3).Value = "鹏华科创主题3年封闭灵活配置混合"
$newQ2.Cells.Item(2, 4).Value = "8.26"
$newQ2.Cells.Item(2, 5).Value = "51.60"
$newQ2.Cells.Item(2, 6).Value = "1.06"
$newQ2.Cells.Item(2, 7).Value = "0.0876"
$newQ2.Cells.Item(2, 8).Value = 7

# Update the "总计" summary sheet: the new 2022-Q4 row is inserted right
# after the header, pushing the old 2022-Q2 row down to row 3.
$summary.Cells.Item(2, 1).Copy()
$summary.Cells.Item(3, 1).PasteSpecial(-4122)

$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(3, 2).Value = "2022-Q2"
$summary.Cells.Item(3, 3).Value = 1
$summary.Cells.Item(3, 4).Value = 0.09

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q4"
$summary.Cells.Item(2, 3).Value = 4
$summary.Cells.Item(2, 4).Value = 0.72
